# Swag Labs login test data: remove the "locked_out_user" row from the
# LoginUsers sheet (row 3), which shifts the rows below it up by one,
# and move the active selection to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing "locked_out_user" / "secret_sauce"
# (row 3) - this shifts rows 4 and 5 up, producing the new A1:B4 table.
$ws.Rows("3:3").Delete()

# Update the active cell / selection shown in the sheet view.
$ws.Range("D12").Select()
